$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 9.916623205318373
$ws.Range("C2").Value = 5.07897920462138
$ws.Range("D2").Value = 14.93946261057526
$ws.Range("E2").Value = 16.35167184383697
$ws.Range("G2").Value = 33.81191038979743
$ws.Range("H2").Value = 15.64677123867274
$ws.Range("I2").Value = 22.45997610142042
$ws.Range("J2").Value = 9.330680799534088
$ws.Range("K2").Value = 10.03002775146818
$ws.Range("O2").Value = 24.47319651619621
$ws.Range("B3").Value = 9.594904838512338
$ws.Range("C3").Value = 4.845612279250438
$ws.Range("D3").Value = 14.87473339026774
$ws.Range("E3").Value = 16.28694576143924
$ws.Range("G3").Value = 33.92431149510465
$ws.Range("H3").Value = 15.70092781821229
$ws.Range("I3").Value = 22.56479203947258
$ws.Range("J3").Value = 9.337108408962457
$ws.Range("K3").Value = 9.808224056255591
$ws.Range("O3").Value = 24.56439798213218
$ws.Range("B4").Value = 9.39281831168061
$ws.Range("C4").Value = 4.696649736166301
$ws.Range("D4").Value = 14.83819637284612
$ws.Range("E4").Value = 16.25074957583737
$ws.Range("G4").Value = 34.0034399210498
$ws.Range("H4").Value = 15.73665458199001
$ws.Range("I4").Value = 22.63340514714515
$ws.Range("J4").Value = 9.342437764196758
$ws.Range("K4").Value = 9.670548136542619
$ws.Range("O4").Value = 24.62548321906054
$ws.Range("B5").Value = 9.30944562326899
$ws.Range("C5").Value = 4.634594608660224
$ws.Range("D5").Value = 14.82412555506296
$ws.Range("E5").Value = 16.23690241849525
$ws.Range("G5").Value = 34.03821771822451
$ws.Range("H5").Value = 15.75183576862053
$ws.Range("I5").Value = 22.6624355331807
$ws.Range("J5").Value = 9.344957576531836
$ws.Range("K5").Value = 9.614144989235394
$ws.Range("O5").Value = 24.65165296675366
$ws.Range("B6").Value = 9.295543805378582
$ws.Range("C6").Value = 4.624211182958214
$ws.Range("D6").Value = 14.82183885689399
$ws.Range("E6").Value = 16.23465796448523
$ws.Range("G6").Value = 34.04414516245772
$ws.Range("H6").Value = 15.75439417599492
$ws.Range("I6").Value = 22.66732062799788
$ws.Range("J6").Value = 9.34539701972189
$ws.Range("K6").Value = 9.604763471851706
$ws.Range("O6").Value = 24.65607549745509
$ws.Range("B7").Value = 9.391697881016286
$ws.Range("C7").Value = 4.695818204321817
$ws.Range("D7").Value = 14.83800328064711
$ws.Range("E7").Value = 16.25055915741232
$ws.Range("G7").Value = 34.00389870811155
$ws.Range("H7").Value = 15.73685680131514
$ws.Range("I7").Value = 22.63379232891306
$ws.Range("J7").Value = 9.342470337601377
$ws.Range("K7").Value = 9.669788573832593
$ws.Range("O7").Value = 24.62583098577174
$ws.Range("B8").Value = 9.806709940076001
$ws.Range("C8").Value = 4.999733694737865
$ws.Range("D8").Value = 14.91648582144156
$ws.Range("E8").Value = 16.32862554243183
$ws.Range("G8").Value = 33.8485608447696
$ws.Range("H8").Value = 15.66493088331392
$ws.Range("I8").Value = 22.49523327602672
$ws.Range("J8").Value = 9.332610274173812
$ws.Range("K8").Value = 9.953902633683493
$ws.Range("O8").Value = 24.50358563399783
$ws.Range("B9").Value = 10.57954597372358
$ws.Range("C9").Value = 5.547920368872359
$ws.Range("D9").Value = 15.0952649069047
$ws.Range("E9").Value = 16.50929511609488
$ws.Range("G9").Value = 33.62465070646476
$ws.Range("H9").Value = 15.54351621238763
$ws.Range("I9").Value = 22.25729277612207
$ws.Range("J9").Value = 9.32422884767689
$ws.Range("K9").Value = 10.49611166325953
$ws.Range("O9").Value = 24.30432400477769
$ws.Range("B10").Value = 11.11651603146042
$ws.Range("C10").Value = 5.918486215111605
$ws.Range("D10").Value = 15.24094712271838
$ws.Range("E10").Value = 16.65803589254189
$ws.Range("G10").Value = 33.5099334485442
$ws.Range("H10").Value = 15.46627884738838
$ws.Range("I10").Value = 22.10307579441757
$ws.Range("J10").Value = 9.324721751895719
$ws.Range("K10").Value = 10.8814184966874
$ws.Range("O10").Value = 24.1827230629198
$ws.Range("B11").Value = 11.35306701689855
$ws.Range("C11").Value = 6.07959224229772
$ws.Range("D11").Value = 15.31013782589491
$ws.Range("E11").Value = 16.72898319803684
$ws.Range("G11").Value = 33.46866829274749
$ws.Range("H11").Value = 15.43373967982484
$ws.Range("I11").Value = 22.03739405910819
$ws.Range("J11").Value = 9.32638256615434
$ws.Range("K11").Value = 11.0531123846319
$ws.Range("O11").Value = 24.1328152252965
$ws.Range("B12").Value = 11.44145880960068
$ws.Range("C12").Value = 6.139493496866756
$ws.Range("D12").Value = 15.33674052035658
$ws.Range("E12").Value = 16.75630367554226
$ws.Range("G12").Value = 33.45462046243342
$ws.Range("H12").Value = 15.42179138058845
$ws.Range("I12").Value = 22.01316569176253
$ws.Range("J12").Value = 9.327217230888097
$ws.Range("K12").Value = 11.11755541289255
$ws.Range("O12").Value = 24.11469646529377
$ws.Range("B13").Value = 11.42247585725872
$ws.Range("C13").Value = 6.12664234257289
$ws.Range("D13").Value = 15.33099358497178
$ws.Range("E13").Value = 16.75039981326657
$ws.Range("G13").Value = 33.45757559544602
$ws.Range("H13").Value = 15.42434804231115
$ws.Range("I13").Value = 22.01835505469179
$ws.Range("J13").Value = 9.327028334341769
$ws.Range("K13").Value = 11.10370285370768
$ws.Range("O13").Value = 24.11856392412906
$ws.Range("B14").Value = 11.36036320919657
$ws.Range("C14").Value = 6.084542688811675
$ws.Range("D14").Value = 15.31231850150907
$ws.Range("E14").Value = 16.73122187187854
$ws.Range("G14").Value = 33.46748089287138
$ws.Range("H14").Value = 15.43274919780854
$ws.Range("I14").Value = 22.0353878678614
$ws.Range("J14").Value = 9.326447116076217
$ws.Range("K14").Value = 11.05842596858027
$ws.Range("O14").Value = 24.13130893275411
$ws.Range("B15").Value = 11.32216102967327
$ws.Range("C15").Value = 6.058610519472845
$ws.Range("D15").Value = 15.30093123480291
$ws.Range("E15").Value = 16.71953344042549
$ws.Range("G15").Value = 33.47375395116649
$ws.Range("H15").Value = 15.43794380513729
$ws.Range("I15").Value = 22.045904833535
$ws.Range("J15").Value = 9.326117871607352
$ws.Range("K15").Value = 11.03061619506196
$ws.Range("O15").Value = 24.13921730225652
$ws.Range("B16").Value = 11.1008949344097
$ws.Range("C16").Value = 5.90780441327622
$ws.Range("D16").Value = 15.23648255800212
$ws.Range("E16").Value = 16.65346388770912
$ws.Range("G16").Value = 33.51285067633727
$ws.Range("H16").Value = 15.46845762392223
$ws.Range("I16").Value = 22.10745830923412
$ws.Range("J16").Value = 9.3246420529545
$ws.Range("K16").Value = 10.87012069263906
$ws.Range("O16").Value = 24.18609371119726
$ws.Range("B17").Value = 10.96312217218323
$ws.Range("C17").Value = 5.813352548852489
$ws.Range("D17").Value = 15.19768040328979
$ws.Range("E17").Value = 16.61376100639956
$ws.Range("G17").Value = 33.5396382044386
$ws.Range("H17").Value = 15.48784203102628
$ws.Range("I17").Value = 22.14636537278788
$ws.Range("J17").Value = 9.324104071819143
$ws.Range("K17").Value = 10.77070107312946
$ws.Range("O17").Value = 24.21623800968539
$ws.Range("B18").Value = 10.88315624257863
$ws.Range("C18").Value = 5.758325128840708
$ws.Range("D18").Value = 15.17563832420375
$ws.Range("E18").Value = 16.5912352330007
$ws.Range("G18").Value = 33.55607312226433
$ws.Range("H18").Value = 15.49923581106061
$ws.Range("I18").Value = 22.1691645957135
$ws.Range("J18").Value = 9.323929913694361
$ws.Range("K18").Value = 10.7131834016896
$ws.Range("O18").Value = 24.23408519474227
$ws.Range("B19").Value = 10.85595945629402
$ws.Range("C19").Value = 5.739574439669045
$ws.Range("D19").Value = 15.16822317931011
$ws.Range("E19").Value = 16.58366221863687
$ws.Range("G19").Value = 33.56181390751075
$ws.Range("H19").Value = 15.50313551595457
$ws.Range("I19").Value = 22.1769562900521
$ws.Range("J19").Value = 9.323894203933344
$ws.Range("K19").Value = 10.69365336506812
$ws.Range("O19").Value = 24.24021528779035
$ws.Range("B20").Value = 10.97786368257189
$ws.Range("C20").Value = 5.823479936176527
$ws.Range("D20").Value = 15.20178253366849
$ws.Range("E20").Value = 16.61795546167762
$ws.Range("G20").Value = 33.5366802290295
$ws.Range("H20").Value = 15.48575323374107
$ws.Range("I20").Value = 22.14218008157892
$ws.Range("J20").Value = 9.324147346137421
$ws.Range("K20").Value = 10.78131947913143
$ws.Range("O20").Value = 24.21297640205574
$ws.Range("B21").Value = 11.37863991264248
$ws.Range("C21").Value = 6.096938628257783
$ws.Range("D21").Value = 15.31779307249544
$ws.Range("E21").Value = 16.73684271683057
$ws.Range("G21").Value = 33.46452856984929
$ws.Range("H21").Value = 15.4302714343382
$ws.Range("I21").Value = 22.03036743667444
$ws.Range("J21").Value = 9.326612257055071
$ws.Range("K21").Value = 11.07174090168862
$ws.Range("O21").Value = 24.12754421970693
$ws.Range("B22").Value = 11.63363483688774
$ws.Range("C22").Value = 6.26920116833556
$ws.Range("D22").Value = 15.39594507463527
$ws.Range("E22").Value = 16.81718136120792
$ws.Range("G22").Value = 33.4265764794726
$ws.Range("H22").Value = 15.39618846673385
$ws.Range("I22").Value = 21.96104482960736
$ws.Range("J22").Value = 9.329422036189687
$ws.Range("K22").Value = 11.25817993996536
$ws.Range("O22").Value = 24.07625836435082
$ws.Range("B23").Value = 11.49819640257139
$ws.Range("C23").Value = 6.177861695951303
$ws.Range("D23").Value = 15.35402661951438
$ws.Range("E23").Value = 16.77406780340419
$ws.Range("G23").Value = 33.44598768937704
$ws.Range("H23").Value = 15.4141798680814
$ws.Range("I23").Value = 21.99769990530669
$ws.Range("J23").Value = 9.327813011491557
$ws.Range("K23").Value = 11.15900027651544
$ws.Range("O23").Value = 24.10321356380751
$ws.Range("B24").Value = 10.97120140261073
$ws.Range("C24").Value = 5.818903602513646
$ws.Range("D24").Value = 15.19992713123759
$ws.Range("E24").Value = 16.61605821281285
$ws.Range("G24").Value = 33.53801430818323
$ws.Range("H24").Value = 15.4866968020464
$ws.Range("I24").Value = 22.14407090912566
$ws.Range("J24").Value = 9.324127360827548
$ws.Range("K24").Value = 10.77652001464104
$ws.Range("O24").Value = 24.21444936514543
$ws.Range("B25").Value = 10.37550461830322
$ws.Range("C25").Value = 5.405085535879191
$ws.Range("D25").Value = 15.04431940957542
$ws.Range("E25").Value = 16.45754612330666
$ws.Range("G25").Value = 33.67652058491199
$ws.Range("H25").Value = 15.57426063273795
$ws.Range("I25").Value = 22.31804561025305
$ws.Range("J25").Value = 9.325325974450042
$ws.Range("K25").Value = 10.35145139510756
$ws.Range("O25").Value = 24.35388428864547
